# Trade #5 closed at 2026-02-17 20:48:03 - unknown UNKNOWN +0.000%
# Also opens a new MarketMaking trade (#66).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet — roll up totals after the newly-closed trade
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1400.32
$wsSummary.Range("B4").Value = 0.11
$wsSummary.Range("B5").Value = 0.07000000000000001
$wsSummary.Range("B6").Value = 33
$wsSummary.Range("B7").Value = 14
$wsSummary.Range("B9").Value = 42.42

# ---------------------------------------------------------------
# Strategy Status sheet — MarketMaking row (row 5)
# ---------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 100.32
$wsStatus.Range("F5").Value = 0.32

# ---------------------------------------------------------------
# All Trades sheet — close out trade #33 (row 34) and append the
# newly opened trade #66 (row 67)
# ---------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")
$wsAll.Range("G34").Value = 0.92
$wsAll.Range("H34").Value = "CLOSED"
$wsAll.Range("I34").Value = 4.5455
$wsAll.Range("J34").Value = 0.04
$wsAll.Range("K34").Value = 100.32
$wsAll.Range("L34").Value = "early_exit"
$wsAll.Range("M34").Value = 0.14

$wsAll.Range("A67").Value = 66
$wsAll.Range("B67").Value = "'2026-02-17"
$wsAll.Range("C67").Value = "20:47:56"
$wsAll.Range("D67").Value = "MarketMaking"
$wsAll.Range("E67").Value = "DOWN"
$wsAll.Range("F67").Value = 0.88
$wsAll.Range("H67").Value = "OPEN"
$wsAll.Range("I67").Value = 0
$wsAll.Range("J67").Value = 0
$wsAll.Range("K67").Value = 100.28
$wsAll.Range("M67").Value = 0
$wsAll.Range("N67").Value = 0
$wsAll.Range("O67").Value = 0
$wsAll.Range("P67").Value = 0.6
$wsAll.Range("Q67").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------
# MarketMaking sheet — append the newly opened trade #66 (row 34)
# ---------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("A34").Value = 66
$wsMM.Range("B34").Value = "'2026-02-17"
$wsMM.Range("C34").Value = "20:47:56"
$wsMM.Range("D34").Value = "MarketMaking"
$wsMM.Range("E34").Value = "DOWN"
$wsMM.Range("F34").Value = 0.88
$wsMM.Range("H34").Value = "OPEN"
$wsMM.Range("I34").Value = 0
$wsMM.Range("J34").Value = 0
$wsMM.Range("K34").Value = 100.28
$wsMM.Range("L34").Value = 0
$wsMM.Range("M34").Value = 0
$wsMM.Range("N34").Value = 0.6
$wsMM.Range("O34").Value = "Normal spread capture: 19600 bps"
$wsMM.Range("Q34").Value = 0
